# [DSD-1088] Trimmed out the data which is exceeding 64 column size on
# mosip-master doc_type table.
#
# The `descr` (column D / C) text for a handful of rows is longer than the
# 64-character column-size limit for this table, so those rows are removed
# entirely:
#   - hin / DOC007  (row 134) "सेवा फोटो आईडी कार्ड जो एक सार्वजनिक उपक्रम ..."
#   - hin / DOC025  (row 147) "राज्य सरकार, सीजीएचएस, ईसीएचएस और ईएसआईसी ..."
#   - tam / DOC007  (row 164) "பொதுத்துறை நிறுவனத்தால் வழங்கப்படும் ..."
#   - tam / DOC016  (row 173) "PSU ஆனது முகவரியுடன் கூடிய சேவை ..."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-up so the earlier row numbers stay valid as we go.
$ws.Rows.Item(173).Delete()
$ws.Rows.Item(164).Delete()
$ws.Rows.Item(147).Delete()
$ws.Rows.Item(134).Delete()

# Re-settle the (cosmetic) row heights that Excel/Calc recomputed for the
# remaining rows once the sheet was edited.
$heights = @{
    2=13.8; 3=13.8; 4=13.8; 5=13.8; 6=13.8; 7=13.8; 8=13.8; 9=13.8; 10=13.8;
    11=13.8; 12=13.8; 13=13.8; 14=13.8; 15=13.8; 16=13.8; 17=13.8; 18=13.8;
    19=13.8; 20=13.8; 21=13.8; 22=13.8; 23=13.8; 24=13.8; 25=13.8; 26=13.8;
    27=13.8; 28=13.8; 29=13.8; 30=13.8; 31=13.8; 32=13.8; 33=13.8; 34=13.8;
    35=13.8; 36=13.8; 37=13.8; 38=13.8; 39=13.8; 40=13.8; 41=13.8; 42=13.8;
    43=13.8; 44=13.8; 45=13.8; 46=13.8; 47=13.8; 48=13.8; 49=13.8; 50=13.8;
    51=13.8; 52=13.8; 53=13.8; 54=13.8; 55=13.8; 56=13.8; 57=13.8; 58=13.8;
    59=13.8; 60=13.8; 61=13.8; 62=13.8; 63=13.8; 64=13.8; 65=13.8; 66=13.8;
    67=13.8; 68=13.8; 69=13.8; 70=16.4; 71=13.8; 72=13.8; 73=13.8; 74=16.4;
    75=13.8; 76=13.8; 77=13.8; 78=13.8; 79=13.8; 80=13.8; 81=13.8; 82=13.8;
    83=16.4; 84=13.8; 85=16.4; 86=16.4; 87=16.4; 88=13.8; 89=13.8; 90=16.4;
    91=13.8; 120=16.4
}

foreach ($r in $heights.Keys) {
    $ws.Rows.Item($r).RowHeight = $heights[$r]
}

# Match the cursor/selection position left behind in the saved file.
$ws.Range("E91").Select()
